$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string with a significant
# trailing zero (e.g. "7.20") must be forced to Text format first, or Excel
# would store them as the number 7.2 and silently drop the trailing zero.
$ws.Range('D2').Value = '66.800.05'
$ws.Range('E2').Value = '  +2.30%  '
$ws.Range('D3').Value = '3.089.90'
$ws.Range('E3').Value = '  +5.40%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '579.78'
$ws.Range('E5').Value = '  +2.12%  '
$ws.Range('D6').Value = '167.98'
$ws.Range('E6').Value = '  +6.13%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.086.28'
$ws.Range('E8').Value = '  +5.35%  '
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('D11').Value = '0.155'
$ws.Range('E11').Value = '  +3.74%  '
$ws.Range('D12').Value = '0.483'
$ws.Range('E12').Value = '  +5.31%  '
$ws.Range('E13').Value = '  +1.85%  '
$ws.Range('D14').Value = '36.49'
$ws.Range('E14').Value = '  +6.65%  '
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('D16').Value = '3.601.02'
$ws.Range('E16').Value = '  +5.12%  '
$ws.Range('D17').Value = '66.803.04'
$ws.Range('E17').Value = '  +2.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.20'
$ws.Range('E18').Value = '  +3.49%  '
$ws.Range('D19').Value = '3.087.52'
$ws.Range('E19').Value = '  +5.12%  '
$ws.Range('D20').Value = '16.15'
$ws.Range('E20').Value = '  +2.89%  '
$ws.Range('D21').Value = '467.13'
$ws.Range('E21').Value = '  +5.13%  '
$ws.Range('D22').Value = '0.714'
$ws.Range('E22').Value = '  +3.83%  '
$ws.Range('D23').Value = '7.52'
$ws.Range('E23').Value = '  +3.80%  '
$ws.Range('D24').Value = '83.89'
$ws.Range('E24').Value = '  +2.19%  '
$ws.Range('E25').Value = '  +6.71%  '
$ws.Range('D26').Value = '13.04'
$ws.Range('E26').Value = '  +7.98%  '
$ws.Range('D27').Value = '10.13'
$ws.Range('E27').Value = '  +0.96%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.00'
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('E30').Value = '  +1.92%  '
$ws.Range('E31').Value = '  +4.01%  '
$ws.Range('E32').Value = '  +1.25%  '
$ws.Range('D33').Value = '28.27'
$ws.Range('E33').Value = '  +4.58%  '
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value = '1.01'
$ws.Range('E36').Value = '  +3.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.90'
$ws.Range('E37').Value = '  +3.35%  '
$ws.Range('D38').Value = '47.27'
$ws.Range('E38').Value = '  +4.65%  '
$ws.Range('E39').Value = '  +6.53%  '
$ws.Range('D40').Value = '0.319'
$ws.Range('E40').Value = '  +6.39%  '
$ws.Range('D41').Value = '50.33'
$ws.Range('E41').Value = '  +1.28%  '
$ws.Range('D42').Value = '0.122'
$ws.Range('E42').Value = '  +1.10%  '
$ws.Range('E43').Value = '  +2.76%  '
$ws.Range('D44').Value = '2.82'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').Value = '0.0361'
$ws.Range('E45').Value = '  +2.74%  '
$ws.Range('D46').Value = '383.93'
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').Value = '2.782.80'
$ws.Range('E47').Value = '  +3.13%  '
$ws.Range('D48').Value = '135.09'
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('D50').Value = '24.98'
$ws.Range('E50').Value = '  +6.95%  '
$ws.Range('E51').Value = '  +2.24%  '
